$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new students were uploaded at the bottom of the roster (rows 38-39),
# each with an ID in column A and a cloudskillsboost public-profile link in
# column B (plain text, like the most recently added rows 35-37 -- no
# hyperlink formatting applied to these trailing rows).
$ws.Range("A38").Value = 220131698
$ws.Range("B38").Value = "https://www.cloudskillsboost.google/public_profiles/425f8493-3bd1-45e9-b20d-0f18af59f453"

$ws.Range("A39").Value = 220106745
$ws.Range("B39").Value = "https://www.cloudskillsboost.google/public_profiles/a03a4b08-27bc-49c4-910c-42251eacaae4"

# Scroll the view down so the newly added rows are visible, and select the
# full used range (matching the author's saved view state).
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("A1:B39").Select()
